# "student and course collection complete"
# - Courses table (rows 17-24) gains a "Prerequisities" field (row 24,
#   replacing the old lone "Year" row) and gets a blue highlight on the
#   "Semester" / new "Prerequisities" field-name cells.
# - Department table (rows 26+) is filled out with Department Head, Office
#   Manager, Professors, Associate Professors, Assistant Professors,
#   Lecturers, Courses and Number of students rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Courses table: add "Prerequisities" as a new attribute row ---
$ws.Range("A24").Value = "Prerequisities"

# Highlight the "Semester" and "Prerequisities" field-name cells with a
# solid blue fill (standard color, RGB 0,176,240 -> 00B0F0).
$ws.Range("A23").Interior.Color = 15773696
$ws.Range("A24").Interior.Color = 15773696

# --- Department table: flesh out the remaining attribute rows ---
$ws.Range("A28").Value = "Department Head"
$ws.Range("A29").Value = "Office Manager"
$ws.Range("A30").Value = "Professors"
$ws.Range("A31").Value = "Associate Professors"
$ws.Range("A32").Value = "Assistant Professors"
$ws.Range("A33").Value = "Lecturers"
$ws.Range("A34").Value = "Courses"
$ws.Range("A35").Value = "Number of students"

# Column A needs to stay wide enough for the new labels.
$ws.Columns("A").ColumnWidth = 17.17

# Leave the cursor where the author left off.
$ws.Range("E17").Select() | Out-Null
